$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("table_id_to_chart")
$ws.Range("G55").Style = "Neutral"
$ws.Range("K55").Style = "Neutral"
